$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 4672.25
$ws.Range("J16").Value = 4630
$ws.Range("L16").Value = 4630
$ws.Range("N16").Value = -5090
$ws.Range("H28").Value = 795
$ws.Range("I28").Value = 472.22223
$ws.Range("K28").Value = 472.22223
$ws.Range("M28").Value = 12.77776999999998
$ws.Range("H64").Value = 3700
$ws.Range("J64").Value = 3933.3333
$ws.Range("L64").Value = 3933.3333
$ws.Range("N64").Value = -4429.3333
$ws.Range("H67").Value = 3700
$ws.Range("J67").Value = 3933.3333
$ws.Range("L67").Value = 3933.3333
$ws.Range("N67").Value = -5649.3333
$ws.Range("H100").Value = 73103.28999999999
$ws.Range("J100").Value = 2375
$ws.Range("L100").Value = 2375
$ws.Range("N100").Value = -3457
$ws.Range("H116").Value = 8004.7144
$ws.Range("I116").Value = 6218.5625
$ws.Range("J116").Value = 10386.25
$ws.Range("K116").Value = 6218.5625
$ws.Range("L116").Value = 10386.25
$ws.Range("M116").Value = -2776.5625
$ws.Range("N116").Value = -17270.25
$ws.Range("H137").Value = 17774.215
$ws.Range("I137").Value = 7607.75
$ws.Range("K137").Value = 22823.25
$ws.Range("M137").Value = -20273.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1130.3846
$ws.Range("I132").Value = 1269.6
$ws.Range("J132").Value = 666.3333
$ws.Range("K132").Value = 3808.8
$ws.Range("L132").Value = 1998.9999
$ws.Range("M132").Value = -1278.8
$ws.Range("N132").Value = -7058.9999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2971.7666
$ws.Range("I105").Value = 3182.889
$ws.Range("K105").Value = 3182.889
$ws.Range("M105").Value = -1435.889
$ws.Range("H135").Value = 39604.684
$ws.Range("J135").Value = 39604.684
$ws.Range("L135").Value = 39604.684
$ws.Range("N135").Value = -49744.684
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 17221.285
$ws.Range("I41").Value = 15099.833
$ws.Range("J41").Value = 29950
$ws.Range("K41").Value = 15099.833
$ws.Range("L41").Value = 29950
$ws.Range("M41").Value = -14671.833
$ws.Range("N41").Value = -30806
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H50").Value = 29999.273
$ws.Range("J50").Value = 29999.273
$ws.Range("L50").Value = 29999.273
$ws.Range("N50").Value = -31249.273
$ws.Range("H59").Value = 60528.227
$ws.Range("J59").Value = 59601.094
$ws.Range("L59").Value = 59601.094
$ws.Range("N59").Value = -61891.094
$ws.Range("H60").Value = 22083.334
$ws.Range("I60").Value = 14166.667
$ws.Range("K60").Value = 14166.667
$ws.Range("M60").Value = -13655.667
$ws.Range("H86").Value = 4995.8887
$ws.Range("I86").Value = 4995.8335
$ws.Range("J86").Value = 4996
$ws.Range("K86").Value = 4995.8335
$ws.Range("L86").Value = 4996
$ws.Range("M86").Value = -3872.8335
$ws.Range("N86").Value = -7242
$ws.Range("H89").Value = 4995.8887
$ws.Range("I89").Value = 4995.8335
$ws.Range("J89").Value = 4996
$ws.Range("K89").Value = 24979.1675
$ws.Range("L89").Value = 24980
$ws.Range("M89").Value = -19363.1675
$ws.Range("N89").Value = -36212
$ws.Range("H99").Value = 10071.027
$ws.Range("I99").Value = 7444.7856
$ws.Range("J99").Value = 11742.272
$ws.Range("K99").Value = 7444.7856
$ws.Range("L99").Value = 11742.272
$ws.Range("M99").Value = -5946.7856
$ws.Range("N99").Value = -14738.272
$ws.Range("H126").Value = 10071.027
$ws.Range("I126").Value = 7444.7856
$ws.Range("J126").Value = 11742.272
$ws.Range("K126").Value = 22334.3568
$ws.Range("L126").Value = 35226.81600000001
$ws.Range("M126").Value = -19864.3568
$ws.Range("N126").Value = -40166.81600000001
$ws.Range("H132").Value = 43027.9
$ws.Range("I132").Value = 38067
$ws.Range("J132").Value = 47058.625
$ws.Range("K132").Value = 114201
$ws.Range("L132").Value = 141175.875
$ws.Range("M132").Value = -111671
$ws.Range("N132").Value = -146235.875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 74344250
$ws.Range("I4").Value = 80456270
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 241368810
$ws.Range("L4").Value = 3000000
$ws.Range("M4").Value = -241368698
$ws.Range("N4").Value = -3000224
$ws.Range("H58").Value = 3500
$ws.Range("I58").Value = 3500
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 10500
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -10372
$ws.Range("N58").ClearContents()
$ws.Range("H104").Value = 2551.5
$ws.Range("I104").Value = 3016
$ws.Range("J104").Value = 229
$ws.Range("K104").Value = 9048
$ws.Range("L104").Value = 687
$ws.Range("M104").Value = -6427
$ws.Range("N104").Value = -5929
$ws.Range("H140").Value = 3738.1667
$ws.Range("J140").Value = 5000
$ws.Range("L140").Value = 15000
$ws.Range("N140").Value = -25360
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 57500
$ws.Range("J104").Value = 57500
$ws.Range("L104").Value = 57500
$ws.Range("N104").Value = -64488
$ws.Range("H113").Value = 2987
$ws.Range("I113").Value = 2974
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2974
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -804
$ws.Range("N113").Value = -7340
$ws.Range("H121").Value = 83093.664
$ws.Range("J121").Value = 83093.664
$ws.Range("L121").Value = 83093.664
$ws.Range("N121").Value = -86587.664
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 6110.5454
$ws.Range("I34").Value = 10500
$ws.Range("J34").Value = 5135.1113
$ws.Range("K34").Value = 10500
$ws.Range("L34").Value = 5135.1113
$ws.Range("M34").Value = -10328
$ws.Range("N34").Value = -5479.1113
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 1685
$ws.Range("I40").Value = 1422.1
$ws.Range("K40").Value = 1422.1
$ws.Range("M40").Value = -1286.1
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 14574.5
$ws.Range("I51").Value = 14574.5
$ws.Range("K51").Value = 14574.5
$ws.Range("M51").Value = -14064.5
$ws.Range("H126").Value = 7010.1
$ws.Range("I126").Value = 7113.9375
$ws.Range("K126").Value = 21341.8125
$ws.Range("M126").Value = -18871.8125
$ws.Range("H132").Value = 25880.8
$ws.Range("I132").Value = 17998.354
$ws.Range("J132").Value = 36188.617
$ws.Range("K132").Value = 53995.062
$ws.Range("L132").Value = 108565.851
$ws.Range("M132").Value = -51465.062
$ws.Range("N132").Value = -113625.851
